# Refresh market-board-derived profit figures (currentAveragePrice* / LevePrice* /
# LeveProfit*) on each job sheet of the Ifrit_Profits workbook, cell by cell, to
# match the latest Universalis pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 1786.35
$ws.Range("I62").Value = 1502.1333
$ws.Range("J62").Value = 2639
$ws.Range("K62").Value = 1502.1333
$ws.Range("L62").Value = 2639
$ws.Range("M62").Value = -878.1333
$ws.Range("N62").Value = -3887
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 1786.35
$ws.Range("I65").Value = 1502.1333
$ws.Range("J65").Value = 2639
$ws.Range("K65").Value = 7510.666499999999
$ws.Range("L65").Value = 13195
$ws.Range("M65").Value = -4390.666499999999
$ws.Range("N65").Value = -19435
# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 406.07693
$ws.Range("I92").Value = 421.1
$ws.Range("J92").Value = 356
$ws.Range("K92").Value = 421.1
$ws.Range("L92").Value = 356
$ws.Range("M92").Value = 826.9
$ws.Range("N92").Value = -2852
# Row 100: Asking for a Friend
$ws.Range("H100").Value = 1404.7142
$ws.Range("I100").Value = 1107.0714
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1107.0714
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -566.0714
$ws.Range("N100").Value = -3082
# Row 101: Edge of the Arcane
$ws.Range("H101").Value = 2316.5715
$ws.Range("I101").Value = 258
$ws.Range("J101").Value = 3140
$ws.Range("K101").Value = 774
$ws.Range("L101").Value = 9420
$ws.Range("M101").Value = 848
$ws.Range("N101").Value = -12664
# Row 125: Body over Mind
$ws.Range("H125").Value = 1886.4
$ws.Range("I125").Value = 366
$ws.Range("J125").Value = 2900
$ws.Range("K125").Value = 3294
$ws.Range("L125").Value = 26100
$ws.Range("M125").Value = -834
$ws.Range("N125").Value = -31020
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 3954.95
$ws.Range("I137").Value = 4817.037
$ws.Range("K137").Value = 14451.111
$ws.Range("M137").Value = -11901.111

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 3678299.2
$ws.Range("I2").Value = 3245
$ws.Range("K2").Value = 3245
$ws.Range("M2").Value = -3132
# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1596.1111
$ws.Range("I102").Value = 1521.3334
$ws.Range("J102").Value = 1970
$ws.Range("K102").Value = 1521.3334
$ws.Range("L102").Value = 1970
$ws.Range("M102").Value = 100.6666
$ws.Range("N102").Value = -5214
# Row 116: No Scope
$ws.Range("H116").Value = 3678299.2
$ws.Range("I116").Value = 3245
$ws.Range("K116").Value = 3245
$ws.Range("M116").Value = -951
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1441.5186
$ws.Range("I122").Value = 872.6667
$ws.Range("J122").Value = 2579.2222
$ws.Range("K122").Value = 2618.0001
$ws.Range("L122").Value = 7737.6666
$ws.Range("M122").Value = -168.0001000000002
$ws.Range("N122").Value = -12637.6666
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 12799922
$ws.Range("I132").Value = 22398086
$ws.Range("J132").Value = 2368.8333
$ws.Range("K132").Value = 67194258
$ws.Range("L132").Value = 7106.499899999999
$ws.Range("M132").Value = -67191728
$ws.Range("N132").Value = -12166.4999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 3678299.2
$ws.Range("I3").Value = 3245
$ws.Range("K3").Value = 3245
$ws.Range("M3").Value = -3131
# Row 94: High Steal
$ws.Range("H94").Value = 1125.2916
$ws.Range("I94").Value = 929.9
$ws.Range("J94").Value = 2102.25
$ws.Range("K94").Value = 929.9
$ws.Range("L94").Value = 2102.25
$ws.Range("M94").Value = -478.9
$ws.Range("N94").Value = -3004.25
# Row 107: The Gold Experience
$ws.Range("H107").Value = 2902.4
$ws.Range("I107").Value = 1913.6666
$ws.Range("J107").Value = 4385.5
$ws.Range("K107").Value = 1913.6666
$ws.Range("L107").Value = 4385.5
$ws.Range("M107").Value = 6.333399999999983
$ws.Range("N107").Value = -8225.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2826213.5
$ws.Range("I31").Value = 1011.0323
$ws.Range("J31").Value = 5954116.5
$ws.Range("K31").Value = 1011.0323
$ws.Range("L31").Value = 5954116.5
$ws.Range("M31").Value = -716.0323
$ws.Range("N31").Value = -5954706.5
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2826213.5
$ws.Range("I34").Value = 1011.0323
$ws.Range("J34").Value = 5954116.5
$ws.Range("K34").Value = 1011.0323
$ws.Range("L34").Value = 5954116.5
$ws.Range("M34").Value = -809.0323
$ws.Range("N34").Value = -5954520.5
# Row 70: A Reward Fitting of the Faithful
$ws.Range("H70").Value = 33500
$ws.Range("J70").Value = 33500
$ws.Range("L70").Value = 33500
$ws.Range("N70").Value = -34130
# Row 73: Just Rewards for Just Devotion (L)
$ws.Range("H73").Value = 33500
$ws.Range("J73").Value = 33500
$ws.Range("L73").Value = 33500
$ws.Range("N73").Value = -35684
# Row 107: Built to Last
$ws.Range("H107").Value = 3418.1304
$ws.Range("I107").Value = 3406.9333
$ws.Range("J107").Value = 3439.125
$ws.Range("K107").Value = 3406.9333
$ws.Range("L107").Value = 3439.125
$ws.Range("M107").Value = -1486.9333
$ws.Range("N107").Value = -7279.125
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3237.6333
$ws.Range("I132").Value = 3076.4
$ws.Range("J132").Value = 3560.1
$ws.Range("K132").Value = 9229.200000000001
$ws.Range("L132").Value = 10680.3
$ws.Range("M132").Value = -6699.200000000001
$ws.Range("N132").Value = -15740.3

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 1157.909
$ws.Range("I68").Value = 859.5789
$ws.Range("J68").Value = 1384.64
$ws.Range("K68").Value = 2578.7367
$ws.Range("L68").Value = 4153.92
$ws.Range("M68").Value = -1767.7367
$ws.Range("N68").Value = -5775.92
# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 1157.909
$ws.Range("I71").Value = 859.5789
$ws.Range("J71").Value = 1384.64
$ws.Range("K71").Value = 7736.2101
$ws.Range("L71").Value = 12461.76
$ws.Range("M71").Value = -3680.2101
$ws.Range("N71").Value = -20573.76

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1135.7142
$ws.Range("I97").Value = 1146.1538
$ws.Range("K97").Value = 1146.1538
$ws.Range("M97").Value = -650.1538
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 4044.25
$ws.Range("I113").Value = 4044.25
$ws.Range("K113").Value = 4044.25
$ws.Range("M113").Value = -1874.25
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 5502.9473
$ws.Range("J122").Value = 3580
$ws.Range("L122").Value = 10740
$ws.Range("N122").Value = -15640
# Row 132: On Board for Lar
$ws.Range("H132").Value = 15386304
$ws.Range("I132").Value = 41668292
$ws.Range("J132").Value = 1725.5366
$ws.Range("K132").Value = 125004876
$ws.Range("L132").Value = 5176.6098
$ws.Range("M132").Value = -125002346
$ws.Range("N132").Value = -10236.6098

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 1958.2
$ws.Range("I61").Value = 1072.75
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 1072.75
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -870.75
$ws.Range("N61").Value = -5904
# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1865.5
$ws.Range("I93").Value = 1452.0588
$ws.Range("K93").Value = 1452.0588
$ws.Range("M93").Value = -204.0588
# Row 113: Peace in Rest
$ws.Range("H113").Value = 1958.2
$ws.Range("I113").Value = 1072.75
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 1072.75
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = 1097.25
$ws.Range("N113").Value = -9840
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 4519.361
$ws.Range("I132").Value = 4900.3096
$ws.Range("J132").Value = 3677.2632
$ws.Range("K132").Value = 14700.9288
$ws.Range("L132").Value = 11031.7896
$ws.Range("M132").Value = -12170.9288
$ws.Range("N132").Value = -16091.7896

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 5889765.5
$ws.Range("I96").Value = 14286629
$ws.Range("J96").Value = 11961
$ws.Range("K96").Value = 14286629
$ws.Range("L96").Value = 11961
$ws.Range("M96").Value = -14285256
$ws.Range("N96").Value = -14707
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 7847758
$ws.Range("I132").Value = 9342201
$ws.Range("J132").Value = 1935.125
$ws.Range("K132").Value = 28026603
$ws.Range("L132").Value = 5805.375
$ws.Range("M132").Value = -28024073
$ws.Range("N132").Value = -10865.375
